$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column contained values like "4-27-2013-14" (game day/month
# concatenated with the season label). That is one calendar day off from
# the real game date because of how the NBA stats site displayed it, so
# replace it with the correct ISO date "2014-04-27" for every data row.
$oldValue = "4-27-2013-14"
$newValue = "2014-04-27"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$lastCol = $used.Column + $used.Columns.Count - 1

# Locate the "Date" header column on row 1 so the fix applies wherever the
# column happens to live.
$dateCol = 0
for ($col = $used.Column; $col -le $lastCol; $col++) {
    $header = $ws.Cells.Item(1, $col).Value2
    if ($header -eq "Date") {
        $dateCol = $col
        break
    }
}

if ($dateCol -gt 0) {
    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, $dateCol)
        if ($cell.Value2 -eq $oldValue) {
            # Prefix with an apostrophe so Excel stores the corrected value
            # as literal text instead of auto-converting the YYYY-MM-DD
            # pattern into a date serial number.
            $cell.Value2 = "'" + $newValue
        }
    }
}
